$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.66"
$ws.Range("E2").Value = "'1.04%"
$ws.Range("D3").Value = "'27.17"
$ws.Range("E3").Value = "'-3.48%"
$ws.Range("D4").Value = "'4.768"
$ws.Range("E4").Value = "'-10.40%"
$ws.Range("D5").Value = "'0.05968"
$ws.Range("E5").Value = "'2.00%"
$ws.Range("D6").Value = "'6.669"
$ws.Range("E6").Value = "'-0.52%"
$ws.Range("D7").Value = "'0.8699"
$ws.Range("E7").Value = "'0.34%"
$ws.Range("D8").Value = "'0.9485"
$ws.Range("E8").Value = "'4.58%"
$ws.Range("D9").Value = "'0.1408"
$ws.Range("E9").Value = "'-0.84%"
$ws.Range("D10").Value = "'0.07160"
$ws.Range("E10").Value = "'0.10%"
$ws.Range("D11").Value = "'0.03173"
$ws.Range("E11").Value = "'-0.22%"
$ws.Range("D12").Value = "'0.09233"
$ws.Range("E12").Value = "'0.13%"
$ws.Range("D13").Value = "'0.001544"
$ws.Range("E13").Value = "'-0.81%"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006081"
$ws.Range("E14").Value = "'-0.07%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006048"
$ws.Range("E15").Value = "'4.04%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.484"
$ws.Range("E16").Value = "'-0.43%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.183"
$ws.Range("E17").Value = "'-1.41%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.240"
$ws.Range("E18").Value = "'1.76%"
$ws.Range("D20").Value = "'0.03686"
$ws.Range("E20").Value = "'6.97%"
$ws.Range("E21").Value = "'-0.66%"
$ws.Range("D22").Value = "'3.804"
$ws.Range("E22").Value = "'7.39%"
$ws.Range("D23").Value = "'0.04224"
$ws.Range("E23").Value = "'1.63%"
$ws.Range("E24").Value = "'0.17%"
$ws.Range("D25").Value = "'0.001221"
$ws.Range("E25").Value = "'-0.53%"
$ws.Range("D26").Value = "'0.004500"
$ws.Range("E26").Value = "'-10.72%"
$ws.Range("E27").Value = "'0.00%"
$ws.Range("E28").Value = "'-22.88%"
$ws.Range("E40").Value = "'-0.98%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1100"
$ws.Range("E41").Value = "'-0.26%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.003987"
$ws.Range("E42").Value = "'-29.68%"
$ws.Range("D43").Value = "'0.002299"
$ws.Range("E43").Value = "'-1.85%"
$ws.Range("D44").Value = "'0.01059"
$ws.Range("E44").Value = "'-2.73%"
$ws.Range("D45").Value = "'0.00005500"
$ws.Range("E45").Value = "'5.34%"
$ws.Range("E46").Value = "'0.11%"
$ws.Range("D47").Value = "'0.08853"
$ws.Range("E47").Value = "'1.21%"
$ws.Range("E48").Value = "'8.69%"
$ws.Range("E49").Value = "'0.11%"
$ws.Range("E50").Value = "'0.11%"
